$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.224.84'
$ws.Range('E2').Value = '  +1.45%  '
$ws.Range('D3').Value = '3.570.23'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '205.59'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +7.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '559.68'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.91%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.606'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.42%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.672'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.80%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '62.98'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +12.45%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.146'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000277'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.74%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.04'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.89%  '
$ws.Range('D14').Value = '4.152.15'
$ws.Range('E14').Value = '  +0.23%  '
$ws.Range('D15').Value = '3.589.57'
$ws.Range('E15').Value = '  +0.47%  '
$ws.Range('B16').Value = 'TRON'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.125'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.10%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.02'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.54%  '
$ws.Range('D18').Value = '68.062.14'
$ws.Range('E18').Value = '  +1.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.09'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.95%  '
$ws.Range('E20').Value = '  -0.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '399.98'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.10'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.18'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.14'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.85'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.36'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.92%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.82'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +4.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.02'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.73%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '720.76'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +12.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.23'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.48'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.04%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.98'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '63.76'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.111'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '40.95'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.42%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.418'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.22'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.10'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +27.28%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '3.157.28'
$ws.Range('E40').Value = '  -1.89%  '
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0723'
$ws.Range('E41').Value = '  -5.41%  '
$ws.Range('E42').Value = '  -1.73%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.56'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.31%  '
$ws.Range('E45').Value = '  +8.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0409'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.69'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.129'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.07'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.04%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '138.70'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.69'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.70%  '
